$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in values for row 4 (D4, E4) and row 22 (D22, E22)
$ws.Range("D4").Value = 5
$ws.Range("E4").Value = 5
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = 5

# Update the active cell / selection to F22
$ws.Range("F22").Select()
